# Cost calculation updated (more decimals)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

# Storage / Pumped-hydro VOM: more decimals
$ws.Range("E8").Value = 0.00000225

# Storage / Li-ion VOM: more decimals (also an order of magnitude smaller)
$ws.Range("E9").Value = 0.00027

# Wind / Offshore CAPEX decreased, shown in red (cost flag)
$ws.Range("C11").Value = 2000
$ws.Range("C11").Font.Color = 255

# PV / Utility: FOM zeroed out, VOM now populated (red flag, shaded)
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.0018
$ws.Range("E13").Font.Color = 255

# PV / Residential CAPEX decreased (red flag); FOM zeroed out, VOM now populated
$ws.Range("C14").Value = 700
$ws.Range("C14").Font.Color = 255
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.0018
$ws.Range("E14").Font.Color = 255

# Update the active selection shown when the workbook is opened
$ws.Range("E18:E19").Select() | Out-Null
